# fix(stok, penjualan): change timeformat
#
# The "penjualan_tanggal" (date) column is removed from the sheet; the
# columns that followed it ("barang_id", "jumlah") shift one slot to the
# left (D, E) and the trailing "jumlah"/"total" data columns (old F, G)
# are cleared out. The G1 header cell keeps its old (bold + centered)
# formatting but loses its text, matching a "select + delete" edit rather
# than a structural column delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# D1 ("penjualan_tanggal") -> "barang_id", reusing E1's current
# (bold, non-centered) formatting.
$ws.Range("E1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "barang_id"

# E1 ("barang_id") -> "jumlah", reusing F1's current (bold,
# non-centered) formatting.
$ws.Range("F1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "jumlah"

# F1 ("jumlah") is no longer used.
$ws.Range("F1").ClearContents()

# G1 ("total") loses its text but keeps its old bold+centered header
# formatting (untouched).
$ws.Range("G1").ClearContents()

$excel.CutCopyMode = 0

# --- Data rows (2-4): shift barang_id/jumlah from E/F into D/E --------
$rows = 2, 3, 4
foreach ($r in $rows) {
    $barangId = $ws.Cells.Item($r, 5).Value2   # old column E = barang_id
    $jumlah = $ws.Cells.Item($r, 6).Value2      # old column F = jumlah

    # Column D: drop the date formatting, reset to the default style,
    # then write the barang_id value that used to live in column E.
    $ws.Cells.Item($r, 4).Style = "Normal"
    $ws.Cells.Item($r, 4).Value = $barangId

    # Column E: jumlah value that used to live in column F.
    $ws.Cells.Item($r, 5).Value = $jumlah

    # Columns F and G no longer hold data.
    $ws.Cells.Item($r, 6).ClearContents()
    $ws.Cells.Item($r, 7).ClearContents()
}
